$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column B slightly (target stored width 14.42578125; COM ColumnWidth
# is quantized to Excel's character/pixel grid, so 13.6666... is the input
# that lands in the nearest achievable bucket, 14.5)
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666

# Update the computed values in A1:B3 (row 4 stays unchanged)
$ws.Range("A1").Value = 0.0054581959047847496
$ws.Range("B1").Value = -0.0054581959500206263

$ws.Range("A2").Value = -0.031817198976671714
$ws.Range("B2").Value = 0.03181719895283288

$ws.Range("A3").Value = 0.065390397539566478
$ws.Range("B3").Value = -0.065390397564560651
